$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1471.3549
$ws.Range("J17").Value = 1234.9667
$ws.Range("L17").Value = 3704.9001
$ws.Range("N17").Value = -4040.9001
$ws.Range("H62").Value = 1450
$ws.Range("I62").Value = 1466.6666
$ws.Range("J62").Value = 1400
$ws.Range("K62").Value = 1466.6666
$ws.Range("L62").Value = 1400
$ws.Range("M62").Value = -842.6666
$ws.Range("N62").Value = -2648
$ws.Range("H65").Value = 1450
$ws.Range("I65").Value = 1466.6666
$ws.Range("J65").Value = 1400
$ws.Range("K65").Value = 7333.333000000001
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = -4213.333000000001
$ws.Range("N65").Value = -13240
$ws.Range("H103").Value = 911.9167
$ws.Range("I103").Value = 844.5
$ws.Range("J103").Value = 1249
$ws.Range("K103").Value = 2533.5
$ws.Range("L103").Value = 3747
$ws.Range("M103").Value = -1947.5
$ws.Range("N103").Value = -4919
$ws.Range("H137").Value = 1786.8334
$ws.Range("I137").Value = 1021
$ws.Range("J137").Value = 2552.6667
$ws.Range("K137").Value = 3063
$ws.Range("L137").Value = 7658.000100000001
$ws.Range("M137").Value = -513
$ws.Range("N137").Value = -12758.0001
$ws.Range("H139").Value = 48899.57
$ws.Range("J139").Value = 48899.57
$ws.Range("L139").Value = 48899.57
$ws.Range("N139").Value = -59179.57
$ws.Range("H141").Value = 934953.75
$ws.Range("I141").Value = 1120840.5
$ws.Range("J141").Value = 5519.6
$ws.Range("K141").Value = 3362521.5
$ws.Range("L141").Value = 16558.8
$ws.Range("M141").Value = -3357341.5
$ws.Range("N141").Value = -26918.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1111610.8
$ws.Range("I2").Value = 1389263.5
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1389263.5
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -1389150.5
$ws.Range("N2").Value = -1226
$ws.Range("H32").Value = 3693.662
$ws.Range("I32").Value = 3008.879
$ws.Range("K32").Value = 3008.879
$ws.Range("M32").Value = -2721.879
$ws.Range("H45").Value = 1641.619
$ws.Range("I45").Value = 1509.8889
$ws.Range("J45").Value = 1740.4166
$ws.Range("K45").Value = 1509.8889
$ws.Range("L45").Value = 1740.4166
$ws.Range("M45").Value = -1132.8889
$ws.Range("N45").Value = -2494.4166
$ws.Range("H61").Value = 1257.6945
$ws.Range("I61").Value = 618.03125
$ws.Range("K61").Value = 618.03125
$ws.Range("M61").Value = -406.03125
$ws.Range("H74").Value = 1227.4681
$ws.Range("I74").Value = 891.9706
$ws.Range("J74").Value = 2104.923
$ws.Range("K74").Value = 891.9706
$ws.Range("L74").Value = 2104.923
$ws.Range("M74").Value = -17.97059999999999
$ws.Range("N74").Value = -3852.923
$ws.Range("H77").Value = 1227.4681
$ws.Range("I77").Value = 891.9706
$ws.Range("J77").Value = 2104.923
$ws.Range("K77").Value = 4459.853
$ws.Range("L77").Value = 10524.615
$ws.Range("M77").Value = -91.85300000000007
$ws.Range("N77").Value = -19260.615
$ws.Range("H102").Value = 1164.25
$ws.Range("I102").Value = 1164.25
$ws.Range("K102").Value = 1164.25
$ws.Range("M102").Value = 457.75
$ws.Range("H116").Value = 1111610.8
$ws.Range("I116").Value = 1389263.5
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1389263.5
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = -1386969.5
$ws.Range("N116").Value = -5588
$ws.Range("H136").Value = 1257.6945
$ws.Range("I136").Value = 618.03125
$ws.Range("K136").Value = 1854.09375
$ws.Range("M136").Value = 695.90625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1111610.8
$ws.Range("I3").Value = 1389263.5
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1389263.5
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -1389149.5
$ws.Range("N3").Value = -1228
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280
$ws.Range("H140").Value = 36666.145
$ws.Range("J140").Value = 41332.6
$ws.Range("L140").Value = 41332.6
$ws.Range("N140").Value = -51692.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 70010
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 70010
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents() | Out-Null
$ws.Range("N23").Value = -70490
$ws.Range("H27").Value = 70010
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 70010
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents() | Out-Null
$ws.Range("N27").Value = -70394
$ws.Range("H31").Value = 2227.762
$ws.Range("J31").Value = 2282.8572
$ws.Range("L31").Value = 2282.8572
$ws.Range("N31").Value = -2872.8572
$ws.Range("H34").Value = 2227.762
$ws.Range("J34").Value = 2282.8572
$ws.Range("L34").Value = 2282.8572
$ws.Range("N34").Value = -2686.8572
$ws.Range("H105").Value = 1075.8572
$ws.Range("I105").Value = 1075.8572
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1075.8572
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents() | Out-Null
$ws.Range("H107").Value = 602.3333
$ws.Range("I107").Value = 602.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 602.3333
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents() | Out-Null
$ws.Range("H134").Value = 2065.8484
$ws.Range("I134").Value = 1990.2609
$ws.Range("K134").Value = 5970.7827
$ws.Range("M134").Value = -3435.7827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 651
$ws.Range("I121").Value = 491.8
$ws.Range("J121").Value = 850
$ws.Range("K121").Value = 1475.4
$ws.Range("L121").Value = 2550
$ws.Range("M121").Value = -165.4000000000001
$ws.Range("N121").Value = -5170
$ws.Range("H126").Value = 4059.8
$ws.Range("I126").Value = 1149.5
$ws.Range("K126").Value = 3448.5
$ws.Range("M126").Value = 1491.5
$ws.Range("H131").Value = 10944.594
$ws.Range("J131").Value = 13877.167
$ws.Range("L131").Value = 41631.501
$ws.Range("N131").Value = -51711.501
$ws.Range("H134").Value = 2151.8215
$ws.Range("I134").Value = 855.3158
$ws.Range("K134").Value = 2565.9474
$ws.Range("M134").Value = 2504.0526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 2100001.8
$ws.Range("I20").Value = 2100001.8
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2100001.8
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents() | Out-Null
$ws.Range("H102").Value = 1872
$ws.Range("I102").Value = 1835.3334
$ws.Range("K102").Value = 1835.3334
$ws.Range("M102").Value = -213.3334
$ws.Range("H107").Value = 533.3333
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 550
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1370
$ws.Range("N107").Value = -4340
$ws.Range("H132").Value = 1427543.5
$ws.Range("I132").Value = 2749402.2
$ws.Range("K132").Value = 8248206.600000001
$ws.Range("M132").Value = -8245676.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2518.9375
$ws.Range("I61").Value = 2164.5
$ws.Range("K61").Value = 2164.5
$ws.Range("M61").Value = -1962.5
$ws.Range("H82").Value = 1868.75
$ws.Range("J82").Value = 2124.5
$ws.Range("L82").Value = 2124.5
$ws.Range("N82").Value = -2846.5
$ws.Range("H85").Value = 1868.75
$ws.Range("J85").Value = 2124.5
$ws.Range("L85").Value = 2124.5
$ws.Range("N85").Value = -4620.5
$ws.Range("H100").Value = 1729.8
$ws.Range("I100").Value = 1583.1666
$ws.Range("K100").Value = 1583.1666
$ws.Range("M100").Value = -1042.1666
$ws.Range("H113").Value = 2518.9375
$ws.Range("I113").Value = 2164.5
$ws.Range("K113").Value = 2164.5
$ws.Range("M113").Value = 5.5
$ws.Range("H132").Value = 1715.5238
$ws.Range("I132").Value = 1118.5588
$ws.Range("J132").Value = 4252.625
$ws.Range("K132").Value = 3355.6764
$ws.Range("L132").Value = 12757.875
$ws.Range("M132").Value = -825.6764000000003
$ws.Range("N132").Value = -17817.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("N21").ClearContents() | Out-Null
$ws.Range("H29").Value = 6725
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 8800
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 8800
$ws.Range("M29").Value = -210
$ws.Range("N29").Value = -9380
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("N35").ClearContents() | Out-Null
$ws.Range("H96").Value = 8527
$ws.Range("I96").Value = 2159.4
$ws.Range("J96").Value = 13833.333
$ws.Range("K96").Value = 2159.4
$ws.Range("L96").Value = 13833.333
$ws.Range("M96").Value = -786.4000000000001
$ws.Range("N96").Value = -16579.333
$ws.Range("H126").Value = 1892.6552
$ws.Range("J126").Value = 4678.4287
$ws.Range("L126").Value = 14035.2861
$ws.Range("N126").Value = -18975.2861
$ws.Range("H132").Value = 1175.6028
$ws.Range("I132").Value = 820.0469000000001
$ws.Range("K132").Value = 2460.1407
$ws.Range("M132").Value = 69.85930000000008
